$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-10-10 Friday" "2025-10-11 Saturday"

Replace-Text "79×94=" "99×59="
Replace-Text "57×62=" "55×34="
Replace-Text "58×92=" "87×32="
Replace-Text "74×60=" "32×56="
Replace-Text "46×92=" "47×73="
Replace-Text "23×67=" "36×61="
Replace-Text "61×71=" "56×14="
Replace-Text "73×50=" "70×19="
Replace-Text "91×45=" "79×88="
Replace-Text "86×71=" "75×56="
Replace-Text "32×18=" "48×13="
Replace-Text "12×48=" "26×70="
Replace-Text "52×48=" "11×14="
Replace-Text "82×39=" "84×93="
Replace-Text "78×33=" "40×65="
Replace-Text "29×74=" "18×22="
Replace-Text "99×19=" "47×76="
Replace-Text "99×75=" "59×76="
Replace-Text "20×30=" "59×28="
Replace-Text "84×25=" "91×50="
Replace-Text "39×36=" "74×81="
Replace-Text "82×56=" "54×63="
Replace-Text "43×48=" "61×97="
Replace-Text "56×72=" "92×15="
Replace-Text "71×22=" "52×17="

Write-Output "Done"
